# Update Name of Algo
# Apply updated RandomForest imputation result values to Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value  = 5.732599999999999
$ws.Range("B3").Value  = 5.94249999999999
$ws.Range("B5").Value  = 4.915100000000004
$ws.Range("E7").Value  = 11.868
$ws.Range("A9").Value  = -20.45599999999997
$ws.Range("E9").Value  = 13.00970000000001
$ws.Range("B11").Value = 5.313399999999997
$ws.Range("B12").Value = 5.466599999999999
$ws.Range("A13").Value = -22.02260000000001
$ws.Range("A16").Value = -20.09529999999999
$ws.Range("A18").Value = -22.67420000000001
$ws.Range("A20").Value = -22.19830000000002
$ws.Range("B21").Value = 5.315499999999999
$ws.Range("E21").Value = 12.60699999999999
